$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.07140851020813
$ws.Range("B1").Value = 2.875567197799683
$ws.Range("C1").Value = 1.993100523948669
$ws.Range("D1").Value = 1.784650087356567
$ws.Range("E1").Value = 1.700886249542236
